$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '72.988.56'
$ws.Range("E2").Value = '  +1.41%  '

# Row 3
$ws.Range("D3").Value = '3.989.85'
$ws.Range("E3").Value = '  -1.01%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.00%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '617.15'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +15.20%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '165.37'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +10.52%  '

# Row 7
$ws.Range("E7").Value = '  -1.40%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.03%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.762'
$ws.Range("D9").Style = "Normal"

# Row 10
$ws.Range("E10").Value = '  -2.24%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '58.59'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +8.49%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000318'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.56%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.25'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.96%  '

# Row 14
$ws.Range("D14").Value = '4.634.20'
$ws.Range("E14").Value = '  -0.68%  '

# Row 15
$ws.Range("D15").Value = '3.995.16'
$ws.Range("E15").Value = '  -1.03%  '

# Row 16
$ws.Range("E16").Value = '  +6.04%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.36'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.12%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '20.69'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.95%  '

# Row 19
$ws.Range("E19").Value = '  -0.02%  '

# Row 20
$ws.Range("D20").Value = '72.921.32'
$ws.Range("E20").Value = '  +1.38%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '442.08'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.44%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.94'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +16.56%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '96.58'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.76%  '

# Row 24
$ws.Range("E24").Value = '  -5.57%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '14.37'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.82%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.16'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.61%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.33'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.66%  '

# Row 28
$ws.Range("B28").Value = 'LEO'
$ws.Range("C28").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.97'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.38%  '

# Row 29
$ws.Range("B29").Value = 'Filecoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.55'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.30%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.23'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.34%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.85'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.56%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '13.90'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.03%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.131'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.39%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '49.13'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.29%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '71.56'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.02%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '643.88'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.76%  '

# Row 37
$ws.Range("D37").Value = '0.0₃0924'
$ws.Range("E37").Value = '  +12.12%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.436'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.34%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.52'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.37%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.148'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.56%  '

# Row 41
$ws.Range("E41").Value = '  -0.14%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.06'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.09%  '

# Row 43
$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.28%  '

# Row 44
$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.28'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.41%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0489'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.08%  '

# Row 46
$ws.Range("E46").Value = '  -0.72%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.96'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +36.39%  '

# Row 48
$ws.Range("E48").Value = '  -0.94%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.41'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.51%  '

# Row 50
$ws.Range("B50").Value = 'Stacks'
$ws.Range("C50").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.09'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.02%  '

# Row 51
$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").Value = '2.902.14'
$ws.Range("E51").Value = '  +0.92%  '
